$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.259
$ws.Range("E2").Value = 0.4425
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 66.8
$ws.Range("L2").Value = 0.4156813939016802
$ws.Range("M2").Value = 9.199999999999999
$ws.Range("N2").Value = 0.008226037195994277
$ws.Range("O2").Value = 0.1377245508982036
$ws.Range("P2").Value = 9.199999999999999
$ws.Range("Q2").Value = 0.008226037195994277
$ws.Range("R2").Value = 0.1377245508982036
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 384
$ws.Range("V2").Value = 0.3433476394849785
$ws.Range("W2").Value = 0.1886611580403794
$ws.Range("X2").Value = 0.1204369700204922
$ws.Range("Y2").Value = 0.06822418801988722
$ws.Range("Z2").Value = 0.07428191071379045
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.07260230968913881
$ws.Range("AC2").Value = -0.07260230968913881
$ws.Range("AD2").Value = 2276.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2276.9
$ws.Range("AG2").Value = 1892.9
$ws.Range("AH2").Value = 0.67060348128295
$ws.Range("AI2").Value = 0.8024882811123251
$ws.Range("AJ2").Value = 0.628598943977684
$ws.Range("AK2").Value = 0.7715729833285778

# Row 3
$ws.Range("D3").Value = 0.394
$ws.Range("E3").Value = 0.627
$ws.Range("K3").Value = 19.7
$ws.Range("L3").Value = 0.652317880794702
$ws.Range("M3").Value = 6.25
$ws.Range("N3").Value = 0.01022745868106693
$ws.Range("O3").Value = 0.317258883248731
$ws.Range("P3").Value = 6.25
$ws.Range("Q3").Value = 0.01022745868106693
$ws.Range("R3").Value = 0.317258883248731
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 37.6
$ws.Range("V3").Value = 0.06152839142529865
$ws.Range("W3").Value = 0.2086864406779661
$ws.Range("X3").Value = 0.06418429986532175
$ws.Range("Y3").Value = 0.1445021408126443
$ws.Range("Z3").Value = 0.1093885830194147
$ws.Range("AB3").Value = 0.06265832702482663
$ws.Range("AC3").Value = -0.06265832702482663
$ws.Range("AD3").Value = 252.1
$ws.Range("AF3").Value = 252.1
$ws.Range("AG3").Value = 214.5
$ws.Range("AH3").Value = 0.2920528266913809
$ws.Range("AI3").Value = 0.5625976344565945
$ws.Range("AJ3").Value = 0.2598110465116279
$ws.Range("AK3").Value = 0.5225334957369062

# Row 4
$ws.Range("D4").Value = 0.124
$ws.Range("E4").Value = 0.258
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 47.1
$ws.Range("L4").Value = 0.3609195402298851
$ws.Range("M4").Value = 2.95
$ws.Range("N4").Value = 0.005815099546619357
$ws.Range("O4").Value = 0.06263269639065817
$ws.Range("P4").Value = 2.95
$ws.Range("Q4").Value = 0.005815099546619357
$ws.Range("R4").Value = 0.06263269639065817
$ws.Range("U4").Value = 346.4
$ws.Range("V4").Value = 0.6828306721860832
$ws.Range("W4").Value = 0.1686358754027927
$ws.Range("X4").Value = 0.1766896401756626
$ws.Range("Y4").Value = -0.008053764772869898
$ws.Range("Z4").Value = 0.06914639961850262
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.08254629235345097
$ws.Range("AC4").Value = -0.08254629235345097
$ws.Range("AD4").Value = 2024.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2024.8
$ws.Range("AG4").Value = 1678.4
$ws.Range("AH4").Value = 0.7996524623830022
$ws.Range("AI4").Value = 0.8474803281433116
$ws.Range("AJ4").Value = 0.7679004437937502
$ws.Range("AK4").Value = 0.8216173878989622

# Remove now-empty cells (AN/AP) for rows 2 and 4
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()